$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the names to their initials
$ws.Range("A2").Value = "JSM"
$ws.Range("A3").Value = "MFBT"

# Move the selection to A4 (was C4)
$ws.Range("A4").Select()
